$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume update (GitHub Actions scheduled refresh).
# For cells whose new text happens to look like a plain number (e.g. "1.005"),
# force the cell to Text format first so Excel keeps it as a string instead of
# auto-converting it to a numeric value (these are formatted price strings, not numbers).

$ws.Range("D2").Value = "24.364.10"
$ws.Range("E2").Value = "  -2.29%  "
$ws.Range("D3").Value = "1.647.10"
$ws.Range("E3").Value = "  -3.73%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.47"
$ws.Range("E5").Value = "  -0.78%  "
$ws.Range("E6").Value = "  +0.42%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3645"
$ws.Range("E7").Value = "  -3.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "46.52"
$ws.Range("E8").Value = "  -6.50%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3225"
$ws.Range("E9").Value = "  -7.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.114"
$ws.Range("E10").Value = "  -8.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06988"
$ws.Range("E11").Value = "  -7.86%  "
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.915"
$ws.Range("E13").Value = "  -7.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.21"
$ws.Range("E14").Value = "  -10.13%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.552"
$ws.Range("E15").Value = "  -7.63%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "1.645.17"
$ws.Range("E16").Value = "  -3.65%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001027"
$ws.Range("E17").Value = "  -9.63%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06570"
$ws.Range("E18").Value = "  -2.38%  "
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "77.60"
$ws.Range("E20").Value = "  -9.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.890"
$ws.Range("E21").Value = "  -8.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.47"
$ws.Range("E22").Value = "  -11.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.45"
$ws.Range("E23").Value = "  -6.72%  "
$ws.Range("D24").Value = "24.365.60"
$ws.Range("E24").Value = "  -2.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.481"
$ws.Range("E25").Value = "  +1.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.284"
$ws.Range("E26").Value = "  -18.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "145.60"
$ws.Range("E27").Value = "  -4.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.43"
$ws.Range("E28").Value = "  -10.48%  "
$ws.Range("D29").Value = "1.827.69"
$ws.Range("E29").Value = "  -3.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "123.14"
$ws.Range("E30").Value = "  -7.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.169"
$ws.Range("E31").Value = "  -6.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.053"
$ws.Range("E32").Value = "  -4.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.602"
$ws.Range("E33").Value = "  -19.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08397"
$ws.Range("E34").Value = "  -4.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.651"
$ws.Range("E35").Value = "  -8.91%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.99"
$ws.Range("E36").Value = "  -14.33%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.256"
$ws.Range("E37").Value = "  -2.52%  "
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.134"
$ws.Range("E38").Value = "  -9.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05950"
$ws.Range("E39").Value = "  -11.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02203"
$ws.Range("E40").Value = "  -9.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2038"
$ws.Range("E41").Value = "  -9.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.049"
$ws.Range("E42").Value = "  -14.13%  "
$ws.Range("E43").Value = "  +0.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5835"
$ws.Range("E44").Value = "  -10.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.752"
$ws.Range("E45").Value = "  -2.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.36"
$ws.Range("E46").Value = "  -11.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5547"
$ws.Range("E47").Value = "  -10.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "121.57"
$ws.Range("E48").Value = "  -6.96%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.927"
$ws.Range("E49").Value = "  -10.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06875"
$ws.Range("E50").Value = "  -6.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.174"
$ws.Range("E51").Value = "  -4.27%  "
